# Append the new Adafruit IO data point as row 96 (A96:F96), mirroring
# the existing rows exactly (plain text cells, no special formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 96

$ws.Cells.Item($newRow, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($newRow, 2).Value = "temperature"

# "25" looks like a number, so a plain .Value assignment would be stored
# as a numeric cell. Route it through a text-formula + paste-values so it
# lands as genuine text without touching the shared style table.
$c3 = $ws.Cells.Item($newRow, 3)
$c3.Formula = '="25"'
$c3.Copy()
$c3.PasteSpecial(-4163)

$ws.Cells.Item($newRow, 4).Value = "N/A"
$ws.Cells.Item($newRow, 5).Value = "N/A"
$ws.Cells.Item($newRow, 6).Value = "N/A"
